# Applies the two changes captured in the commit's XML diff:
#  1. Re-caches the "datetimeFigureOut" date field text, on the slide
#     master and on every slide layout, from 9/24/2020 -> 11/16/2020.
#  2. Fixes the label on slide 1's "RDS DB" callout from
#     "instance standby (multi-AZ)" to "instance replica (multi-AZ)".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder ("datetimeFigureOut" field) on master + layouts
# ---------------------------------------------------------------------
$oldDate = "9/24/2020"
$newDate = "11/16/2020"
$ppPlaceholderDate = 16

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)

        $placeholderType = -1
        try { $placeholderType = $shp.PlaceholderFormat.Type } catch { }

        if ($placeholderType -eq $ppPlaceholderDate -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $full = $tr.Characters(1, $tr.Length)
                $full.Text = $newDate
            }
        }
    }
}

# Slide master
Update-DateShapes $p.SlideMaster.Shapes

# Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShapes $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------
# 2) "instance standby (multi-AZ)" -> "instance replica (multi-AZ)"
# ---------------------------------------------------------------------
$oldLabel = "instance standby (multi-AZ)"
$newLabel = "instance replica (multi-AZ)"

$slide = $p.Slides.Item(1)
$group = $slide.Shapes.Item("Group 22")
$textBox = $group.GroupItems.Item("TextBox 19")
$range = $textBox.TextFrame.TextRange

$startPos = $range.Text.IndexOf($oldLabel) + 1
if ($startPos -gt 0) {
    $target = $range.Characters($startPos, $oldLabel.Length)
    $target.Text = $newLabel
}
